$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: column A = document type, column B = description
$data = @(
    @("Types of Testing Documents", "Description"),
    @("Test policy", "It is a high-level document which describes principles, methods and all the important testing goals of the organization."),
    @("Test strategy", "A high-level document which identifies the Test Levels (types) to be executed for the project."),
    @("Test plan", "A test plan is a complete planning document which contains the scope, approach, resources, schedule, etc. of testing activities."),
    @("Requirements Traceability Matrix", "This is a document which connects the requirements to the test cases."),
    @("Test Scenario", "Test scenario is an item or event of a software system which could be verified by one or more Test cases."),
    @("Test case", "It is a group of input values, execution preconditions, expected execution postconditions and results. It is developed for a Test Scenario."),
    @("Test Data", "Test Data is a data which exists before a test is executed. It used to execute the test case."),
    @("Defect Report", "Defect report is a documented report of any flaw in a Software System which fails to perform its expected function."),
    @("Test summary report", "Test summary report is a high-level document which summarizes testing activities conducted as well as the test result.")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Apply border to the whole used range A1:B10
$rng = $ws.Range("A1:B10")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2

# Apply yellow fill to header row
$header = $ws.Range("A1:B1")
$header.Interior.Color = 65535

# Column widths
$ws.Columns.Item(1).ColumnWidth = 28.33203125
$ws.Columns.Item(2).ColumnWidth = 113.44140625

# Selection
$ws.Range("A6").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
